$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the floating point value on A10 (tiny precision correction)
$ws.Range("A10").Value = 45863.70854388889

# Append the new row 11 with the latest sensor reading
$ws.Range("A11").Value = 45863.75018738768
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("B11").Value = 2025
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 18.72
$ws.Range("E11").Value = 73.69
$ws.Range("F11").Value = 12.3
$ws.Range("G11").Value = 6.3
$ws.Range("H11").Value = "ESE"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "18:00:16"
